# Fixed a bug in weightBranch
# The bug caused rows 2-23 (the per-symbol weight table, columns A-F) to be
# written out of order. This corrects the data so that each row contains the
# weight-branch tuple that actually belongs to it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1001, 18, 30, 75, 60, 72),
    @(501,  9,  52, 30, 75, 45),
    @(801,  3,  67, 65, 52, 45),
    @(1203, 3,  15, 15, 15, 15),
    @(901,  16, 15, 45, 60, 60),
    @(301,  6,  45, 30, 60, 45),
    @(201,  9,  30, 15, 45, 30),
    @(1202, 2,  10, 10, 10, 10),
    @(101,  9,  30, 15, 60, 15),
    @(902,  1,  0,  0,  0,  0),
    @(601,  9,  60, 67, 60, 42),
    @(401,  9,  48, 67, 75, 45),
    @(701,  3,  90, 45, 97, 15),
    @(1201, 2,  10, 10, 10, 10),
    @(1101, 0,  15, 30, 30, 0),
    @(1,    0,  2,  2,  2,  2),
    @(502,  0,  4,  0,  0,  0),
    @(802,  0,  4,  5,  4,  0),
    @(3,    0,  3,  3,  3,  3),
    @(2,    0,  2,  2,  2,  2),
    @(602,  0,  0,  4,  0,  9),
    @(402,  0,  0,  4,  0,  0)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
